$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: "Writing" / Mar 20 session - update the time period text
$ws.Range("C17").Value = "14:00-16:00               Mar 20"

# Row 18: "Coding" / Mar 21 session - update the time period text
$ws.Range("C18").Value = "12:00-16:00                  Mar 21"

# Row 19: "Meeting" / Mar 22 session - update the time period text and
# give it the time number format (h:mm), matching the edited workbook.
$ws.Range("C19").Value = "14:00-15:00            Mar 22"
$ws.Range("C19").NumberFormat = "h:mm"

# Row 20: "Coding " / Mar 22 session - update the time period text
$ws.Range("C20").Value = "17:00-20:00                 Mar 22"

# Row 21: was an unused template row, now filled in with a new log entry
$ws.Range("B21").Value = "Coding"
$ws.Range("C21").Value = "14:00-18:00                Mar 24"
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = "Home"
$ws.Range("F21").Value = "Continuing gridding code"
$ws.Range("G21").Value = "Progess on code "

# Move the active selection to C21, matching the saved cursor position
$ws.Range("C21").Select()
